$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1: highlight a few rows with green fill (style used elsewhere: fillId 2 -> RGB 00B050) ---
$greenCells = @("B14", "B21", "B22", "B23")
foreach ($addr in $greenCells) {
    $ws1.Range($addr).Interior.Color = 0x50B000
}

# --- Hoja2: populate with new content ---
$ws2.Range("A1").Value = "tablas  que considero ue faltan"

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "acredores"

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "secciones"

$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "departamentos"

$ws2.Range("A6").Value = 4
$ws2.Range("A7").Value = 5
$ws2.Range("A8").Value = 6
$ws2.Range("A9").Value = 7
$ws2.Range("A10").Value = 8
$ws2.Range("A11").Value = 9

$ws2.Columns.Item(2).ColumnWidth = 14.5703125

# --- Selections / active sheet ---
$ws1.Range("A26").Select()
$ws2.Range("C6").Select()
$ws2.Activate()

$wb.Save()
